$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description text for DE1 (Cus_Nme)
$ws.Range("C2").Value = "Customer Name  and will be used for verification, deposit , deposit cheque, get cheque, withdraw cheque, withdraw money, transfer money, review transactions and for all activities "

# Update description text for DE2 (Acc_num)
$ws.Range("C3").Value = "Account number of the customer will be used to deposit cheque, get cheque book, deposit money, withdrawal, review transactions, verififcation, and other banking activities. "

# Update description text for DE9 (Bill_type)
$ws.Range("C10").Value = "This utility specifies the type of bills that needed for bill payment process"

# Update description text for DE15 (Loan_purp)
$ws.Range("C16").Value = "This signifies the purpose of apply loan like cbuying new car or house"

# Add new row 18 for DE17 (Bill_amt)
$ws.Range("A18").Value = "DE17"
$ws.Range("B18").Value = "Bill_amt"
$ws.Range("C18").Value = "This is used to specify bill amount"

# Update the view selection/scroll position: select C2 so topLeftCell resets and selection moves to C2
$ws.Range("C2").Select()
